$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 29 (this shifts the previous rows 29-33 down to 32-36)
$ws.Rows("29:31").Insert()

# New data block (weekly update), dated 44460 (2021-09-21)
$newRows = @(
    @{ Row=29; H="Sin especificar"; I="Primera"; J=250; K=1900; L=1900; M=1900; P=1900 },
    @{ Row=30; H="Sin especificar"; I="Segunda"; J=160; K=1700; L=1700; M=1700; P=1700 },
    @{ Row=31; H="Sin especificar"; I="Tercera"; J=106; K=1500; L=1500; M=1500; P=1500 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 9
    $ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = 44460
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 300000000
    $ws.Cells.Item($row, 7).Value = "Espárragos"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/kilo"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

Write-Host "done"
